# Learning Journal: Week 2 -> Week 1 content rewrite
$d = $word.ActiveDocument

# 1. Heading: "Week 2" -> "Week 1"
$d.Content.Find.Execute("Week 2", $false, $false, $false, $false, $false, $true, 1, $false, "Week 1", 2)

# 2. "I started watching videos..." -> "I revised all the knowledge..."
$d.Content.Find.Execute(
    "I started watching videos about how to deal with CMS Joomla on website.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "I revised all the knowledge how to code HTML files, css files and javascript by watching videos and doing hand-on practicals that exist on the web.",
    2)

# 3. "I takes several activities..." (spans 3 runs incl. proofErr) -> single sentence
$d.Content.Find.Execute(
    "I takes several activities during the working process, doing practical is one of them. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Since I did the course web design, I can go back and go through all the slide and lectures.",
    2)

# 4. "I found out that I can import..." -> brother anecdote
$d.Content.Find.Execute(
    "I found out that I can import some free Joomla templates from the outside webpage.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "I have a brother that his major is web design, so I ask him how to make the site easy to see, proper functionality and easy to navigate.",
    2)

# 5. "It took me more than 2 days..." -> "It took me around 2 to 2.5 hours..."
$d.Content.Find.Execute(
    "It took me more than 2 days to have my about us and contact us pages done.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "It took me around 2 to 2.5 hours to get the revision completed.",
    2)

# 6. Insert a brand-new paragraph right after the revision-time sentence, about Siteground.
$rng = $d.Content
$rng.Find.Execute("It took me around 2 to 2.5 hours to get the revision completed.")
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(7)
$newRange = $newPara.Range
$newRange.InsertAfter("Since Siteground cannot be reused, I have to find another way to publish my web pages.")

# 7. Big "It's difficult to get familiar..." paragraph (spans 4 runs) -> single reflective paragraph
$d.Content.Find.Execute(
    "It’s difficult to get familiar with the tools in Joomla. With the week 1 experience, I did the work as that video was playing, again, that way keeps me engaged and get more familiar with CMS. Because Joomla is going to be done for the assignment 1 as well.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "I found bored watching videos. So I just practise and do the work from the video that I just watch rather than siting for several hours to watching lecture. That way keeps me engaged and gain my muscle memory to do some works.",
    2)

# 8. Remove the unused "No Spacing" style definition from styles.xml
$style = $d.Styles.Item("No Spacing")
if ($style -ne $null) {
    $style.Delete()
}
